# Auto-generated Excel COM-interop edit script
# Applies meteocat daily summary refresh: updated extraction timestamps
# and the handful of measurement values that shifted between runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 18:47:41"
$ws.Range("E3").Value = "2026-02-06 18:47:44"
$ws.Range("E4").Value = "2026-02-06 18:47:46"
$ws.Range("J4").Value = "997.3 hPa"
$ws.Range("E5").Value = "2026-02-06 18:47:49"
$ws.Range("J5").Value = "997.5 hPa"
$ws.Range("E6").Value = "2026-02-06 18:47:51"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "48%"
$ws.Range("J6").Value = "998.7 hPa"
$ws.Range("E7").Value = "2026-02-06 18:47:54"
$ws.Range("E8").Value = "2026-02-06 18:47:56"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "76%"
$ws.Range("O8").Value = "10.3 °C"
$ws.Range("E9").Value = "2026-02-06 18:47:58"
$ws.Range("E10").Value = "2026-02-06 18:48:01"
$ws.Range("E11").Value = "2026-02-06 18:48:03"
$ws.Range("J11").Value = "998.7 hPa"
$ws.Range("E12").Value = "2026-02-06 18:48:06"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "59%"
$ws.Range("O12").Value = "14.1 °C"
$ws.Range("E13").Value = "2026-02-06 18:48:08"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "76%"
$ws.Range("E14").Value = "2026-02-06 18:48:10"
$ws.Range("E15").Value = "2026-02-06 18:48:13"
$ws.Range("J15").Value = "997.7 hPa"
$ws.Range("E16").Value = "2026-02-06 18:48:15"
$ws.Range("E17").Value = "2026-02-06 18:48:18"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "82%"
$ws.Range("J17").Value = "998.9 hPa"
$ws.Range("E18").Value = "2026-02-06 18:48:20"
$ws.Range("N18").Value = "-6.6 °C 18:00 TU"
$ws.Range("O18").Value = "-4.5 °C"
$ws.Range("E19").Value = "2026-02-06 18:48:23"
$ws.Range("E20").Value = "2026-02-06 18:48:25"
$ws.Range("E21").Value = "2026-02-06 18:48:27"
$ws.Range("J21").Value = "997.9 hPa"
$ws.Range("E22").Value = "2026-02-06 18:48:30"
$ws.Range("K22").Value = "11.5 MJ/m2"
$ws.Range("O22").Value = "10.7 °C"
$ws.Range("E23").Value = "2026-02-06 18:48:32"
$ws.Range("J23").Value = "997.7 hPa"
$ws.Range("E24").Value = "2026-02-06 18:48:34"
$ws.Range("J24").Value = "997.1 hPa"
$ws.Range("E25").Value = "2026-02-06 18:48:37"
$ws.Range("J25").Value = "998.4 hPa"
$ws.Range("E26").Value = "2026-02-06 18:48:39"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "81%"
$ws.Range("O26").Value = "-0.9 °C"
$ws.Range("E27").Value = "2026-02-06 18:48:41"
$ws.Range("J27").Value = "997.7 hPa"
$ws.Range("E28").Value = "2026-02-06 18:48:44"
$ws.Range("J28").Value = "999.8 hPa"
$ws.Range("O28").Value = "5.2 °C"
$ws.Range("E29").Value = "2026-02-06 18:48:46"
$ws.Range("E30").Value = "2026-02-06 18:48:49"
$ws.Range("L30").Value = "49.3 km/h - 316º 18:20 TU"
$ws.Range("E31").Value = "2026-02-06 18:48:51"
$ws.Range("E32").Value = "2026-02-06 18:48:53"
$ws.Range("J32").Value = "999.0 hPa"
$ws.Range("O32").Value = "15.6 °C"
$ws.Range("E33").Value = "2026-02-06 18:48:56"
$ws.Range("E34").Value = "2026-02-06 18:48:58"
$ws.Range("E35").Value = "2026-02-06 18:49:01"
$ws.Range("E36").Value = "2026-02-06 18:49:03"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "48%"
$ws.Range("O36").Value = "13.0 °C"
